# "Fruta / hortaliza, semanal"
# The weekly refresh reshuffles the existing price records (columns D, L, M,
# N, O, P, Q, S, T) across rows 2-26 into a new order; the per-row "key"
# columns (A, B, C, E-K, R) are unaffected, only the variable observation
# columns move between rows. Write each moved value explicitly per cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44498
$ws.Cells.Item(2, 12).Value = "Segunda"
$ws.Cells.Item(2, 13).Value = 300
$ws.Cells.Item(2, 14).Value = 19000
$ws.Cells.Item(2, 15).Value = 20000
$ws.Cells.Item(2, 16).Value = 19500
$ws.Cells.Item(2, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(2, 19).Value = 2438
$ws.Cells.Item(2, 20).Value = 8

$ws.Cells.Item(3, 4).Value = 44895
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 22000
$ws.Cells.Item(3, 15).Value = 22500
$ws.Cells.Item(3, 16).Value = 22250
$ws.Cells.Item(3, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(3, 19).Value = 2781
$ws.Cells.Item(3, 20).Value = 8

$ws.Cells.Item(4, 4).Value = 44880
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22500
$ws.Cells.Item(4, 16).Value = 22250
$ws.Cells.Item(4, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(4, 19).Value = 2781
$ws.Cells.Item(4, 20).Value = 8

$ws.Cells.Item(5, 4).Value = 44159
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 300
$ws.Cells.Item(5, 14).Value = 2000
$ws.Cells.Item(5, 15).Value = 2100
$ws.Cells.Item(5, 16).Value = 2050
$ws.Cells.Item(5, 17).Value = "`$/kilo (en caja de 14 kilos)"
$ws.Cells.Item(5, 19).Value = 2050
$ws.Cells.Item(5, 20).Value = 1

$ws.Cells.Item(6, 4).Value = 44876
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 22000
$ws.Cells.Item(6, 15).Value = 22500
$ws.Cells.Item(6, 16).Value = 22250
$ws.Cells.Item(6, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(6, 19).Value = 2781
$ws.Cells.Item(6, 20).Value = 8

$ws.Cells.Item(7, 4).Value = 44516
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 18000
$ws.Cells.Item(7, 15).Value = 19000
$ws.Cells.Item(7, 16).Value = 18500
$ws.Cells.Item(7, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(7, 19).Value = 2312
$ws.Cells.Item(7, 20).Value = 8

$ws.Cells.Item(8, 4).Value = 44873
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 300
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 22500
$ws.Cells.Item(8, 16).Value = 22250
$ws.Cells.Item(8, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(8, 19).Value = 2781
$ws.Cells.Item(8, 20).Value = 8

$ws.Cells.Item(9, 4).Value = 44509
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 19000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 19500
$ws.Cells.Item(9, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(9, 19).Value = 2438
$ws.Cells.Item(9, 20).Value = 8

$ws.Cells.Item(10, 4).Value = 44890
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 200
$ws.Cells.Item(10, 14).Value = 22000
$ws.Cells.Item(10, 15).Value = 22500
$ws.Cells.Item(10, 16).Value = 22250
$ws.Cells.Item(10, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(10, 19).Value = 2781
$ws.Cells.Item(10, 20).Value = 8

$ws.Cells.Item(11, 4).Value = 44530
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 200
$ws.Cells.Item(11, 14).Value = 19000
$ws.Cells.Item(11, 15).Value = 20000
$ws.Cells.Item(11, 16).Value = 19500
$ws.Cells.Item(11, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(11, 19).Value = 2438
$ws.Cells.Item(11, 20).Value = 8

$ws.Cells.Item(12, 4).Value = 44530
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 16000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 16000
$ws.Cells.Item(12, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(12, 19).Value = 2000
$ws.Cells.Item(12, 20).Value = 8

$ws.Cells.Item(13, 4).Value = 44495
$ws.Cells.Item(13, 12).Value = "Segunda"
$ws.Cells.Item(13, 13).Value = 270
$ws.Cells.Item(13, 14).Value = 19000
$ws.Cells.Item(13, 15).Value = 20000
$ws.Cells.Item(13, 16).Value = 19556
$ws.Cells.Item(13, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(13, 19).Value = 2444
$ws.Cells.Item(13, 20).Value = 8

$ws.Cells.Item(14, 4).Value = 44162
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 2000
$ws.Cells.Item(14, 15).Value = 2100
$ws.Cells.Item(14, 16).Value = 2050
$ws.Cells.Item(14, 17).Value = "`$/kilo (en caja de 14 kilos)"
$ws.Cells.Item(14, 19).Value = 2050
$ws.Cells.Item(14, 20).Value = 1

$ws.Cells.Item(15, 4).Value = 44533
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 300
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 19000
$ws.Cells.Item(15, 16).Value = 18500
$ws.Cells.Item(15, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(15, 19).Value = 2312
$ws.Cells.Item(15, 20).Value = 8

$ws.Cells.Item(16, 4).Value = 44533
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 100
$ws.Cells.Item(16, 14).Value = 16000
$ws.Cells.Item(16, 15).Value = 16000
$ws.Cells.Item(16, 16).Value = 16000
$ws.Cells.Item(16, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(16, 19).Value = 2000
$ws.Cells.Item(16, 20).Value = 8

$ws.Cells.Item(17, 4).Value = 44519
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 400
$ws.Cells.Item(17, 14).Value = 21000
$ws.Cells.Item(17, 15).Value = 22000
$ws.Cells.Item(17, 16).Value = 21500
$ws.Cells.Item(17, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(17, 19).Value = 2688
$ws.Cells.Item(17, 20).Value = 8

$ws.Cells.Item(18, 4).Value = 44519
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 18000
$ws.Cells.Item(18, 15).Value = 18000
$ws.Cells.Item(18, 16).Value = 18000
$ws.Cells.Item(18, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(18, 19).Value = 2250
$ws.Cells.Item(18, 20).Value = 8

$ws.Cells.Item(19, 4).Value = 44523
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 400
$ws.Cells.Item(19, 14).Value = 21000
$ws.Cells.Item(19, 15).Value = 22000
$ws.Cells.Item(19, 16).Value = 21500
$ws.Cells.Item(19, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(19, 19).Value = 2688
$ws.Cells.Item(19, 20).Value = 8

$ws.Cells.Item(20, 4).Value = 44523
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(20, 14).Value = 18000
$ws.Cells.Item(20, 15).Value = 18000
$ws.Cells.Item(20, 16).Value = 18000
$ws.Cells.Item(20, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(20, 19).Value = 2250
$ws.Cells.Item(20, 20).Value = 8

$ws.Cells.Item(21, 4).Value = 44512
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 300
$ws.Cells.Item(21, 14).Value = 19000
$ws.Cells.Item(21, 15).Value = 20000
$ws.Cells.Item(21, 16).Value = 19500
$ws.Cells.Item(21, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(21, 19).Value = 2438
$ws.Cells.Item(21, 20).Value = 8

$ws.Cells.Item(22, 4).Value = 44488
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 160
$ws.Cells.Item(22, 14).Value = 17000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 17500
$ws.Cells.Item(22, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(22, 19).Value = 2188
$ws.Cells.Item(22, 20).Value = 8

$ws.Cells.Item(23, 4).Value = 44491
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 200
$ws.Cells.Item(23, 14).Value = 18000
$ws.Cells.Item(23, 15).Value = 19000
$ws.Cells.Item(23, 16).Value = 18500
$ws.Cells.Item(23, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(23, 19).Value = 2312
$ws.Cells.Item(23, 20).Value = 8

$ws.Cells.Item(24, 4).Value = 44894
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 200
$ws.Cells.Item(24, 14).Value = 22000
$ws.Cells.Item(24, 15).Value = 22500
$ws.Cells.Item(24, 16).Value = 22250
$ws.Cells.Item(24, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(24, 19).Value = 2781
$ws.Cells.Item(24, 20).Value = 8

$ws.Cells.Item(25, 4).Value = 44526
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 300
$ws.Cells.Item(25, 14).Value = 21000
$ws.Cells.Item(25, 15).Value = 21000
$ws.Cells.Item(25, 16).Value = 21000
$ws.Cells.Item(25, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(25, 19).Value = 2625
$ws.Cells.Item(25, 20).Value = 8

$ws.Cells.Item(26, 4).Value = 44505
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 300
$ws.Cells.Item(26, 14).Value = 19000
$ws.Cells.Item(26, 15).Value = 20000
$ws.Cells.Item(26, 16).Value = 19500
$ws.Cells.Item(26, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(26, 19).Value = 2438
$ws.Cells.Item(26, 20).Value = 8
